$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 3.7
$ws.Range("G2").Value = 3.8
$ws.Range("H2").Value = 2.04
$ws.Range("N2").Value = 3.8
$ws.Range("P2").Value = 1.98
$ws.Range("Q2").Value = 1.94
$ws.Range("R2").Value = 1.37
$ws.Range("S2").Value = 3.45
$ws.Range("T2").Value = 1.81
$ws.Range("U2").Value = 2.12
$ws.Range("Z2").Value = 13.5
$ws.Range("AA2").Value = 27
$ws.Range("AB2").Value = 14
$ws.Range("AJ2").Value = 80
$ws.Range("AK2").Value = 48
$ws.Range("AN2").Value = 48
$ws.Range("AO2").Value = 18.5

# Row 3
$ws.Range("F3").Value = 2.54
$ws.Range("G3").Value = 2.64
$ws.Range("I3").Value = 2.78
$ws.Range("V3").Value = 1.56
$ws.Range("W3").Value = 1.6
$ws.Range("AG3").Value = 13
$ws.Range("AO3").Value = 16.5

# Row 5
$ws.Range("F5").Value = 4.7
$ws.Range("G5").Value = 8.4
$ws.Range("H5").Value = 1.9
$ws.Range("I5").Value = 2.46
$ws.Range("J5").Value = 2.54
$ws.Range("K5").Value = 3.85

# Row 6
$ws.Range("K6").Value = 4.1
$ws.Range("Q6").Value = 1.98

# Row 9
$ws.Range("F9").Value = 3.05
$ws.Range("H9").Value = 2.82
$ws.Range("I9").Value = 3.05
$ws.Range("Q9").Value = 2.78

# Row 10
$ws.Range("I10").Value = 3.45

# Row 11
$ws.Range("Q11").Value = 1.71
$ws.Range("AJ11").Value = 12.5

# Row 12
$ws.Range("AH12").Value = 27

# Row 13
$ws.Range("O13").Value = 1.55
$ws.Range("P13").Value = 1.57
$ws.Range("Q13").Value = 2.68
$ws.Range("U13").Value = 1.75
$ws.Range("X13").Value = 8.6
$ws.Range("AE13").Value = 90
$ws.Range("AH13").Value = 25
$ws.Range("AM13").Value = 240
